$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021023025749022
$ws.Range("D2").Value = 1.026857640674566
$ws.Range("E2").Value = 1.021951966342955
$ws.Range("F2").Value = 1.032116099227337
$ws.Range("I2").Value = 1.030361160355375
$ws.Range("J2").Value = 1.026216834019222
$ws.Range("K2").Value = 1.029679062523881
$ws.Range("L2").Value = 1.024787780392354
$ws.Range("M2").Value = 1.034922259029289
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022035308878848
$ws.Range("D3").Value = 1.027620776174032
$ws.Range("E3").Value = 1.022812020828123
$ws.Range("F3").Value = 1.033349849340407
$ws.Range("I3").Value = 1.030586067626785
$ws.Range("J3").Value = 1.026866229871619
$ws.Range("K3").Value = 1.030249921204315
$ws.Range("L3").Value = 1.025454251091353
$ws.Range("M3").Value = 1.035963576800368
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022690289298438
$ws.Range("D4").Value = 1.028114135043577
$ws.Range("E4").Value = 1.023368899236907
$ws.Range("F4").Value = 1.034148103976442
$ws.Range("I4").Value = 1.030729709575228
$ws.Range("J4").Value = 1.027285843640974
$ws.Range("K4").Value = 1.030618207240159
$ws.Range("L4").Value = 1.025885232049177
$ws.Range("M4").Value = 1.036636760968425
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022965634922377
$ws.Range("D5").Value = 1.028321436907811
$ws.Range("E5").Value = 1.023603097918481
$ws.Range("F5").Value = 1.034483675871122
$ws.Range("I5").Value = 1.030789644338343
$ws.Range("J5").Value = 1.027462108071867
$ws.Range("K5").Value = 1.030772771469811
$ws.Range("L5").Value = 1.026066351428872
$ws.Range("M5").Value = 1.036919620342932
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.023011866189614
$ws.Range("D6").Value = 1.028356237555486
$ws.Range("E6").Value = 1.023642425993547
$ws.Range("F6").Value = 1.034540019056433
$ws.Range("I6").Value = 1.030799681116232
$ws.Range("J6").Value = 1.027491695358308
$ws.Range("K6").Value = 1.030798708036184
$ws.Range("L6").Value = 1.026096758347196
$ws.Range("M6").Value = 1.036967105062466
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022693968512018
$ws.Range("D7").Value = 1.028116905439403
$ws.Range("E7").Value = 1.023372028269995
$ws.Range("F7").Value = 1.034152587959674
$ws.Range("I7").Value = 1.030730512204129
$ws.Range("J7").Value = 1.027288199450057
$ws.Range("K7").Value = 1.030620273569789
$ws.Range("L7").Value = 1.025887652431097
$ws.Range("M7").Value = 1.036640541126946
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.021365138770246
$ws.Range("D8").Value = 1.027115636766294
$ws.Range("E8").Value = 1.022242549751743
$ws.Range("F8").Value = 1.032533064755057
$ws.Range("I8").Value = 1.030437559585888
$ws.Range("J8").Value = 1.026436422151552
$ws.Range("K8").Value = 1.029872213919613
$ws.Range("L8").Value = 1.025013073018771
$ws.Range("M8").Value = 1.035274305651307
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019023301711725
$ws.Range("D9").Value = 1.025347926024502
$ws.Range("E9").Value = 1.020255087863558
$ws.Range("F9").Value = 1.029678720373751
$ws.Range("I9").Value = 1.029906891755308
$ws.Range("J9").Value = 1.024930982611179
$ws.Range("K9").Value = 1.028545648330748
$ws.Range("L9").Value = 1.02346989778703
$ws.Range("M9").Value = 1.032862057811752
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.017461885053878
$ws.Range("D10").Value = 1.02416723791518
$ws.Range("E10").Value = 1.018932040513843
$ws.Range("F10").Value = 1.02777539502124
$ws.Range("I10").Value = 1.029543414756537
$ws.Range("J10").Value = 1.023924342570065
$ws.Range("K10").Value = 1.027655653204197
$ws.Range("L10").Value = 1.022439750198952
$ws.Range("M10").Value = 1.031250638827235
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016785725396323
$ws.Range("D11").Value = 1.023655468594212
$ws.Range("E11").Value = 1.018359608842898
$ws.Range("F11").Value = 1.026951116252102
$ws.Range("I11").Value = 1.029383727610241
$ws.Range("J11").Value = 1.023487743101465
$ws.Range("K11").Value = 1.027268946873311
$ws.Range("L11").Value = 1.021993364302187
$ws.Range("M11").Value = 1.030552092907545
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016534560770622
$ws.Range("D12").Value = 1.023465296480981
$ws.Range("E12").Value = 1.018147051172517
$ws.Range("F12").Value = 1.026644921926504
$ws.Range("I12").Value = 1.02932406740875
$ws.Range("J12").Value = 1.023325462826981
$ws.Range("K12").Value = 1.027125106940796
$ws.Range("L12").Value = 1.021827508000152
$ws.Range("M12").Value = 1.030292501877016
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016588436841444
$ws.Range("D13").Value = 1.023506092576274
$ws.Range("E13").Value = 1.018192642389191
$ws.Range("F13").Value = 1.026710602600262
$ws.Range("I13").Value = 1.029336880347109
$ws.Range("J13").Value = 1.023360277374097
$ws.Range("K13").Value = 1.027155970140797
$ws.Range("L13").Value = 1.021863086941554
$ws.Range("M13").Value = 1.030348190481914
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016764964232427
$ws.Range("D14").Value = 1.023639750500649
$ws.Range("E14").Value = 1.018342037351225
$ws.Range("F14").Value = 1.026925806560371
$ws.Range("I14").Value = 1.029378803122271
$ws.Range("J14").Value = 1.023474331162115
$ws.Range("K14").Value = 1.027257061097967
$ws.Range("L14").Value = 1.021979655561222
$ws.Range("M14").Value = 1.0305306374861
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016873727402466
$ws.Range("D15").Value = 1.02372209121038
$ws.Range("E15").Value = 1.018434093663435
$ws.Range("F15").Value = 1.027058398062462
$ws.Range("I15").Value = 1.029404587378264
$ws.Range("J15").Value = 1.023544589182519
$ws.Range("K15").Value = 1.027319320079858
$ws.Range("L15").Value = 1.022051470885671
$ws.Range("M15").Value = 1.030643033174633
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017506758329226
$ws.Range("D16").Value = 1.024201191367889
$ws.Range("E16").Value = 1.018970040567531
$ws.Range("F16").Value = 1.027830096898063
$ws.Range("I16").Value = 1.029553964237325
$ws.Range("J16").Value = 1.023953303120903
$ws.Range("K16").Value = 1.027681289573074
$ws.Range("L16").Value = 1.022469368522424
$ws.Range("M16").Value = 1.031296982258285
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017903826634185
$ws.Range("D17").Value = 1.024501578390591
$ws.Range("E17").Value = 1.019306348382218
$ws.Range("F17").Value = 1.028314128303832
$ws.Range("I17").Value = 1.029647048935288
$ws.Range("J17").Value = 1.024209486455568
$ws.Range("K17").Value = 1.027907986864687
$ws.Range("L17").Value = 1.022731417714858
$ws.Range("M17").Value = 1.031706974691625
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018135424489402
$ws.Range("D18").Value = 1.024676738510691
$ws.Range("E18").Value = 1.019502555069406
$ws.Range("F18").Value = 1.028596443546143
$ws.Range("I18").Value = 1.029701121815175
$ws.Range("J18").Value = 1.024358844612065
$ws.Range("K18").Value = 1.028040086971638
$ws.Range("L18").Value = 1.022884235188965
$ws.Range("M18").Value = 1.031946040112206
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018214392474965
$ws.Range("D19").Value = 1.02473645499711
$ws.Range("E19").Value = 1.01956946393671
$ws.Range("F19").Value = 1.028692703743671
$ws.Range("I19").Value = 1.029719521635837
$ws.Range("J19").Value = 1.024409760136888
$ws.Range("K19").Value = 1.028085107876799
$ws.Range("L19").Value = 1.022936336670112
$ws.Range("M19").Value = 1.032027542347841
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017861225521101
$ws.Range("D20").Value = 1.024469354930871
$ws.Range("E20").Value = 1.019270261179776
$ws.Range("F20").Value = 1.02826219757572
$ws.Range("I20").Value = 1.029637084773963
$ws.Range("J20").Value = 1.024182007570587
$ws.Range("K20").Value = 1.02788367768948
$ws.Range("L20").Value = 1.022703305561753
$ws.Range("M20").Value = 1.031662994276506
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016712981562543
$ws.Range("D21").Value = 1.023600393728328
$ws.Range("E21").Value = 1.01829804235265
$ws.Range("F21").Value = 1.02686243493571
$ws.Range("I21").Value = 1.029366467440311
$ws.Range("J21").Value = 1.02344074813962
$ws.Range("K21").Value = 1.027227297849487
$ws.Range("L21").Value = 1.021945330351417
$ws.Range("M21").Value = 1.03047691471189
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.015990983756229
$ws.Range("D22").Value = 1.023053591106815
$ws.Range("E22").Value = 1.017687168820079
$ws.Range("F22").Value = 1.025982229308188
$ws.Range("I22").Value = 1.029194321935673
$ws.Range("J22").Value = 1.022974065274339
$ws.Range("K22").Value = 1.026813448855533
$ws.Range("L22").Value = 1.021468479534299
$ws.Range("M22").Value = 1.029730484710149
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.016373733453141
$ws.Range("D23").Value = 1.02334350423649
$ws.Range("E23").Value = 1.018010966500164
$ws.Range("F23").Value = 1.026448854779622
$ws.Range("I23").Value = 1.029285768831634
$ws.Range("J23").Value = 1.023221521795741
$ws.Range("K23").Value = 1.027032947630211
$ws.Range("L23").Value = 1.021721293829054
$ws.Range("M23").Value = 1.030126247614299
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017880475146721
$ws.Range("D24").Value = 1.0244839154814
$ws.Range("E24").Value = 1.019286567297095
$ws.Range("F24").Value = 1.028285662873288
$ws.Range("I24").Value = 1.029641587834977
$ws.Range("J24").Value = 1.024194424309744
$ws.Range("K24").Value = 1.027894662356605
$ws.Range("L24").Value = 1.022716008330556
$ws.Range("M24").Value = 1.031682867368113
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019628755742264
$ws.Range("D25").Value = 1.025805314529045
$ws.Range("E25").Value = 1.020768556798239
$ws.Range("F25").Value = 1.030416707462433
$ws.Range("I25").Value = 1.030045792639809
$ws.Range("J25").Value = 1.025320706801653
$ws.Range("K25").Value = 1.028889589147865
$ws.Range("L25").Value = 1.023869087249479
$ws.Range("M25").Value = 1.033486251570873
